$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.007.48"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").Value = "2.429.01"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.62"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.96"
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.511"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.19"
$ws.Range("E10").Value = "  +2.89%  "
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("E12").Value = "  +2.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.63"
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("E14").Value = "  +1.99%  "
$ws.Range("D15").Value = "2.801.91"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("D16").Value = "2.490.92"
$ws.Range("E16").Value = "  +3.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.830"
$ws.Range("E17").Value = "  +2.36%  "
$ws.Range("D18").Value = "44.040.67"
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.20"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").Value = "0.0₃0903"
$ws.Range("E21").Value = "  +1.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.31"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.42"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("E24").Value = "  +2.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.47"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.10"
$ws.Range("E27").Value = "  +1.21%  "
$ws.Range("E28").Value = "  -6.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.46"
$ws.Range("E29").Value = "  +3.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.68"
$ws.Range("E30").Value = "  +3.53%  "
$ws.Range("E31").Value = "  +17.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.62"
$ws.Range("E32").Value = "  +8.06%  "
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0759"
$ws.Range("E35").Value = "  +4.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.92"
$ws.Range("E36").Value = "  +3.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "131.67"
$ws.Range("E37").Value = "  +22.38%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.45"
$ws.Range("E38").Value = "  +1.48%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.93"
$ws.Range("E39").Value = "  +4.31%  "
$ws.Range("E40").Value = "  -1.11%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.55"
$ws.Range("E42").Value = "  -4.73%  "
$ws.Range("D43").Value = "1.952.04"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0286"
$ws.Range("E44").Value = "  +2.01%  "
$ws.Range("E45").Value = "  +2.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.87"
$ws.Range("E46").Value = "  +4.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.32"
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("D48").Value = "2.661.36"
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("E49").Value = "  +6.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.10"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.91"
$ws.Range("E51").Value = "  +0.85%  "
